$d = $word.ActiveDocument

# 1. Replace the long "1268 participants ... Other ... 38.54 (" run sequence
#    with the simplified sentence. We search the whole paragraph range
#    (the FirstParagraph style paragraph) and replace text up through "38.54 (".
$d.Content.Find.Execute(
    "1268 participants were recruited through Amazon Mechanical Turk. A total of 242 participants were excluded before analyses based on the same criteria as Study 1: 14 were excluded because they did not indicate they were American or lived in the United States, 5 were excluded for indicating “Other” for their gender, 155 were excluded for using a phone or tablet to complete the survey, and 68 were excluded for an incomplete survey. The final sample consisted of 1026 participants (50.58% women), with an average age of 38.54 (",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Participants were recruited on Amazon Mechanical Turk using the same pre-screening criteria as Study 1. The final sample consisted of 1026 participants (50.58% women), with an average age of 38.54 (",
    2
)

# 2. "comprehension check questions" -> "comprehension questions"
$d.Content.Find.Execute(
    "comprehension check questions", $true, $false, $false, $false, $false,
    $true, 1, $false, "comprehension questions", 2
)

# 3. "one round per times table" -> "one round per multiplication table"
$d.Content.Find.Execute(
    "one round per times table", $true, $false, $false, $false, $false,
    $true, 1, $false, "one round per multiplication table", 2
)
